$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Lampadas": update existing rows 2-6 and append rows 7-8
# ---------------------------------------------------------------
$wsLamp = $wb.Worksheets.Item("Lampadas")

$wsLamp.Range("A2").Value = "jo"
$wsLamp.Range("B2").Value = "lampada3"
$wsLamp.Range("E2").Value = 48
$wsLamp.Range("H2").Value = "Azul"

$wsLamp.Range("A3").Value = "ji"
$wsLamp.Range("B3").Value = "pi"
$wsLamp.Range("E3").Value = 59
$wsLamp.Range("H3").Value = "Vermelho"

$wsLamp.Range("A4").Value = "ji"
$wsLamp.Range("B4").Value = "po"
$wsLamp.Range("E4").Value = 100
$wsLamp.Range("H4").Value = "Branco"

$wsLamp.Range("A5").Value = "jp"
$wsLamp.Range("B5").Value = "dsakodas"
$wsLamp.Range("E5").Value = 50
$wsLamp.Range("H5").Value = "Branco"

$wsLamp.Range("A6").Value = "jp"
$wsLamp.Range("B6").Value = "vvcx"
$wsLamp.Range("E6").Value = 78
$wsLamp.Range("H6").Value = "Verde"

$wsLamp.Range("A7").Value = "jo"
$wsLamp.Range("B7").Value = "Bom dia"
$wsLamp.Range("E7").Value = 0
$wsLamp.Range("H7").Value = "branca"

$wsLamp.Range("A8").Value = "ji"
$wsLamp.Range("B8").Value = "alo"
$wsLamp.Range("E8").Value = 0
$wsLamp.Range("H8").Value = "branca"

# ---------------------------------------------------------------
# Sheet "Ares Condicionados": update rows 2-3, delete row 4
# ---------------------------------------------------------------
$wsAr = $wb.Worksheets.Item("Ares Condicionados")

$wsAr.Range("A2").Value = "jo"
$wsAr.Range("B2").Value = "vai dar"
$wsAr.Range("D2").Value = 14
$wsAr.Range("E2").Value = 65

$wsAr.Range("A3").Value = "ji"
$wsAr.Range("B3").Value = "que bom"

$wsAr.Rows.Item(4).Delete()
